# updt alert ref; logic if stor/ANF not used
#
# Trim the parenthetical "(MSIX app packages)" / "(User Profiles/FSLogix)"
# qualifiers off the four storage-alert Description cells (F8:F11) so the
# wording is generic regardless of whether ANF or Azure Files is in use.
# Suggestive Action text (G column) is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F8").Value  = "Storage for the follow Azure NetApp volume is critically low"
$ws.Range("F10").Value = "Storage for the follow Azure Files share is critically low"
$ws.Range("F11").Value = "Storage for the follow Azure Files share is moderately low"
$ws.Range("F9").Value  = "Storage for the follow Azure NetApp volume is moderately low"

# Move the saved selection/scroll position of the sheet from the bottom
# (G23) back up to the storage section (F10), matching the author's view
# when they made this edit.
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
$ws.Range("F10").Select()
